$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.288.99"
$ws.Range("E2").Value = "  -3.08%  "
$ws.Range("D3").Value = "3.139.30"
$ws.Range("E3").Value = "  -4.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'525.53"
$ws.Range("E5").Value = "  -5.53%  "
$ws.Range("D6").Value = "'134.65"
$ws.Range("E6").Value = "  -4.79%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.138.57"
$ws.Range("E8").Value = "  -4.92%  "
$ws.Range("D9").Value = "'0.442"
$ws.Range("E9").Value = "  -5.18%  "
$ws.Range("E10").Value = "  -8.17%  "
$ws.Range("E11").Value = "  -8.49%  "
$ws.Range("D12").Value = "'0.379"
$ws.Range("E12").Value = "  -6.82%  "
$ws.Range("D13").Value = "3.674.80"
$ws.Range("E13").Value = "  -4.84%  "
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "'25.53"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "3.132.44"
$ws.Range("E16").Value = "  -5.18%  "
$ws.Range("D17").Value = "58.229.04"
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "'0.0000153"
$ws.Range("E18").Value = "  -6.99%  "
$ws.Range("D19").Value = "'5.78"
$ws.Range("E19").Value = "  -5.12%  "
$ws.Range("D20").Value = "'13.05"
$ws.Range("E20").Value = "  -6.28%  "
$ws.Range("D21").Value = "'7.93"
$ws.Range("E21").Value = "  -7.38%  "
$ws.Range("D22").Value = "'343.84"
$ws.Range("E22").Value = "  -7.86%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'0.506"
$ws.Range("E24").Value = "  -4.71%  "
$ws.Range("D25").Value = "'67.73"
$ws.Range("E25").Value = "  -8.46%  "
$ws.Range("D26").Value = "3.259.82"
$ws.Range("E26").Value = "  -5.14%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.168"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0960"
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'6.81"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -8.08%  "
$ws.Range("D33").Value = "'6.86"
$ws.Range("E33").Value = "  -8.32%  "
$ws.Range("D34").Value = "'21.35"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").Value = "'1.22"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "'4.80"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").Value = "'156.69"
$ws.Range("E37").Value = "  -5.59%  "
$ws.Range("D38").Value = "'6.25"
$ws.Range("E38").Value = "  -6.11%  "
$ws.Range("D39").Value = "'1.37"
$ws.Range("E39").Value = "  -9.58%  "
$ws.Range("D40").Value = "'0.0686"
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").Value = "3.168.98"
$ws.Range("E41").Value = "  -4.72%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'24.37"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'40.41"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").Value = "'0.692"
$ws.Range("E44").Value = "  -7.74%  "
$ws.Range("D45").Value = "'1.08"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "'3.90"
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'1.44"
$ws.Range("E48").Value = "  -8.01%  "
$ws.Range("D49").Value = "2.264.97"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").Value = "'6.22"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").Value = "'20.69"
$ws.Range("E51").Value = "  -2.64%  "

Write-Output "Applied 100 cell updates"
